$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 20 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A20").Value = "Kan je 500 liter antivries bestellen?"
$logs.Range("B20").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C20").Value = "Hoi Johan, `nZou je 500 liter antivries kunnen bestellen voor de wintervoorraad?`nDank je wel!`nGroeten, `nMarc `nSent using {0}"
$logs.Range("D20").Value = "Bestelling / Levering"
$logs.Range("E20").Value = "Bedankt voor je bericht. Ik neem dit z.s.m. in behandeling."
$logs.Range("F20").Value = "2025-06-26 21:29:25"
$logs.Range("G20").Value = "Ja"
$logs.Range("H20").Value = "Ja"
$logs.Range("I20").Value = "Nee"

# --- Extend conditional formatting ranges to include the new row 20 ---
$fc = $logs.Range("D2:D19").FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($logs.Range("D2:D20"))
}

$fc = $logs.Range("G2:G19").FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($logs.Range("G2:G20"))
}

$fc = $logs.Range("H2:H19").FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($logs.Range("H2:H20"))
}

$fc = $logs.Range("I2:I19").FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($logs.Range("I2:I20"))
}

# --- Dashboard sheet: increment "Bestelling / Levering" count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 15
